$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74-177 down to 75-178.
$ws.Rows(74).Insert()

# Populate the newly inserted row 74 with the new data point.
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 44546
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = 100112017
$ws.Cells.Item(74, 7).Value = "Apio"
$ws.Cells.Item(74, 8).Value = "Americana (o)"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 450
$ws.Cells.Item(74, 11).Value = 7500
$ws.Cells.Item(74, 12).Value = 8000
$ws.Cells.Item(74, 13).Value = 7722
$ws.Cells.Item(74, 14).Value = "$/docena de matas"
$ws.Cells.Item(74, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(74, 16).Value = 1287
$ws.Cells.Item(74, 17).Value = 6
$ws.Cells.Item(74, 18).Value = "Hortaliza"
